$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C holds a "Förändrad" (last changed) date serial number.
# Every data row (2 through 149) gets bumped from 46061 (2026-02-08)
# to 46062 (2026-02-09).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 149 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
